# Scheduled runner: refresh computed profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across all job-leve sheets after a market-board
# price sync. Only numeric value cells change; layout/formatting untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 4134.3706
$ws.Range("I41").Value = 316.54544
$ws.Range("J41").Value = 6759.125
$ws.Range("K41").Value = 316.54544
$ws.Range("L41").Value = 6759.125
$ws.Range("M41").Value = 123.45456
$ws.Range("N41").Value = -7639.125
$ws.Range("H58").Value = 1661.4286
$ws.Range("I58").Value = 1210
$ws.Range("K58").Value = 3630
$ws.Range("M58").Value = -3480
$ws.Range("H86").Value = 5836.8
$ws.Range("I86").Value = 3416.111
$ws.Range("J86").Value = 7817.364
$ws.Range("K86").Value = 3416.111
$ws.Range("L86").Value = 7817.364
$ws.Range("M86").Value = -2293.111
$ws.Range("N86").Value = -10063.364
$ws.Range("H89").Value = 5836.8
$ws.Range("I89").Value = 3416.111
$ws.Range("J89").Value = 7817.364
$ws.Range("K89").Value = 17080.555
$ws.Range("L89").Value = 39086.82
$ws.Range("M89").Value = -11464.555
$ws.Range("N89").Value = -50318.82
$ws.Range("H125").Value = 50000268
$ws.Range("I125").Value = 309.14285
$ws.Range("J125").Value = 166666830
$ws.Range("K125").Value = 2782.28565
$ws.Range("L125").Value = 1500001470
$ws.Range("M125").Value = -322.2856500000003
$ws.Range("N125").Value = -1500006390

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1388.091
$ws.Range("I45").Value = 1213.1666
$ws.Range("J45").Value = 1598
$ws.Range("K45").Value = 1213.1666
$ws.Range("L45").Value = 1598
$ws.Range("M45").Value = -836.1666
$ws.Range("N45").Value = -2352

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 839.5217
$ws.Range("I64").Value = 777.7692
$ws.Range("J64").Value = 919.8
$ws.Range("K64").Value = 777.7692
$ws.Range("L64").Value = 919.8
$ws.Range("M64").Value = -552.7692
$ws.Range("N64").Value = -1369.8
$ws.Range("H67").Value = 839.5217
$ws.Range("I67").Value = 777.7692
$ws.Range("J67").Value = 919.8
$ws.Range("K67").Value = 777.7692
$ws.Range("L67").Value = 919.8
$ws.Range("M67").Value = 2.230800000000045
$ws.Range("N67").Value = -2479.8
$ws.Range("H99").Value = 3732.111
$ws.Range("I99").Value = 5099
$ws.Range("J99").Value = 998.3333
$ws.Range("K99").Value = 5099
$ws.Range("L99").Value = 998.3333
$ws.Range("M99").Value = -3601
$ws.Range("N99").Value = -3994.3333
$ws.Range("H134").Value = 2658.6558
$ws.Range("I134").Value = 1073.5814
$ws.Range("J134").Value = 6445.222
$ws.Range("K134").Value = 3220.7442
$ws.Range("L134").Value = 19335.666
$ws.Range("M134").Value = -685.7442000000001
$ws.Range("N134").Value = -24405.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2519.7273
$ws.Range("I31").Value = 1597.2572
$ws.Range("J31").Value = 3561.2258
$ws.Range("K31").Value = 1597.2572
$ws.Range("L31").Value = 3561.2258
$ws.Range("M31").Value = -1302.2572
$ws.Range("N31").Value = -4151.2258
$ws.Range("H34").Value = 2519.7273
$ws.Range("I34").Value = 1597.2572
$ws.Range("J34").Value = 3561.2258
$ws.Range("K34").Value = 1597.2572
$ws.Range("L34").Value = 3561.2258
$ws.Range("M34").Value = -1395.2572
$ws.Range("N34").Value = -3965.2258
$ws.Range("H62").Value = 4632852
$ws.Range("I62").Value = 9261706
$ws.Range("J62").Value = 3998
$ws.Range("K62").Value = 9261706
$ws.Range("L62").Value = 3998
$ws.Range("M62").Value = -9261082
$ws.Range("N62").Value = -5246
$ws.Range("H65").Value = 4632852
$ws.Range("I65").Value = 9261706
$ws.Range("J65").Value = 3998
$ws.Range("K65").Value = 46308530
$ws.Range("L65").Value = 19990
$ws.Range("M65").Value = -46305410
$ws.Range("N65").Value = -26230
$ws.Range("H68").Value = 26787.6
$ws.Range("J68").Value = 26787.6
$ws.Range("L68").Value = 26787.6
$ws.Range("N68").Value = -28285.6
$ws.Range("H71").Value = 26787.6
$ws.Range("J71").Value = 26787.6
$ws.Range("L71").Value = 80362.79999999999
$ws.Range("N71").Value = -87850.79999999999
$ws.Range("H74").Value = 13969.333
$ws.Range("J74").Value = 13969.333
$ws.Range("L74").Value = 13969.333
$ws.Range("N74").Value = -15717.333
$ws.Range("H77").Value = 13969.333
$ws.Range("J77").Value = 13969.333
$ws.Range("L77").Value = 41907.999
$ws.Range("N77").Value = -50643.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 640.6429000000001
$ws.Range("J5").Value = 635
$ws.Range("L5").Value = 1905
$ws.Range("N5").Value = -2129
$ws.Range("H38").Value = 62.5
$ws.Range("I38").Value = 82.666664
$ws.Range("J38").Value = 42.333332
$ws.Range("K38").Value = 247.999992
$ws.Range("L38").Value = 126.999996
$ws.Range("M38").Value = 99.00000800000001
$ws.Range("N38").Value = -820.999996
$ws.Range("H88").Value = 4474.9165
$ws.Range("J88").Value = 4474.9165
$ws.Range("L88").Value = 13424.7495
$ws.Range("N88").Value = -14280.7495
$ws.Range("H91").Value = 4474.9165
$ws.Range("J91").Value = 4474.9165
$ws.Range("L91").Value = 13424.7495
$ws.Range("N91").Value = -16388.7495
$ws.Range("H92").Value = 772.1539
$ws.Range("I92").Value = 675
$ws.Range("J92").Value = 789.8182
$ws.Range("K92").Value = 2025
$ws.Range("L92").Value = 2369.4546
$ws.Range("M92").Value = -777
$ws.Range("N92").Value = -4865.4546
$ws.Range("H97").Value = 553.9
$ws.Range("I97").Value = 209.16667
$ws.Range("J97").Value = 1071
$ws.Range("K97").Value = 627.50001
$ws.Range("L97").Value = 3213
$ws.Range("M97").Value = -131.50001
$ws.Range("N97").Value = -4205
$ws.Range("H122").Value = 206.14285
$ws.Range("I122").Value = 204.4
$ws.Range("J122").Value = 210.5
$ws.Range("K122").Value = 1839.6
$ws.Range("L122").Value = 1894.5
$ws.Range("M122").Value = 610.3999999999999
$ws.Range("N122").Value = -6794.5
$ws.Range("H131").Value = 961.97
$ws.Range("I131").Value = 1916.6666
$ws.Range("J131").Value = 901.0319
$ws.Range("K131").Value = 5749.9998
$ws.Range("L131").Value = 2703.0957
$ws.Range("M131").Value = -709.9997999999996
$ws.Range("N131").Value = -12783.0957
$ws.Range("H132").Value = 4299.9443
$ws.Range("I132").Value = 866.55554
$ws.Range("J132").Value = 7733.3335
$ws.Range("K132").Value = 7798.99986
$ws.Range("L132").Value = 69600.0015
$ws.Range("M132").Value = -5268.99986
$ws.Range("N132").Value = -74660.0015
$ws.Range("H133").Value = 5449.7144
$ws.Range("I133").Value = 4786.6665
$ws.Range("J133").Value = 5947
$ws.Range("K133").Value = 14359.9995
$ws.Range("L133").Value = 17841
$ws.Range("M133").Value = -9299.999500000002
$ws.Range("N133").Value = -27961
$ws.Range("H135").Value = 640.6429000000001
$ws.Range("J135").Value = 635
$ws.Range("L135").Value = 5715
$ws.Range("N135").Value = -10785

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4651.069
$ws.Range("I132").Value = 5437.2666
$ws.Range("J132").Value = 3808.7144
$ws.Range("K132").Value = 16311.7998
$ws.Range("L132").Value = 11426.1432
$ws.Range("M132").Value = -13781.7998
$ws.Range("N132").Value = -16486.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1908.3334
$ws.Range("I61").Value = 1902.9412
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1902.9412
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1700.9412
$ws.Range("N61").Value = -2404
$ws.Range("H68").Value = 2286.8
$ws.Range("I68").Value = 1922.5
$ws.Range("J68").Value = 2529.6667
$ws.Range("K68").Value = 1922.5
$ws.Range("L68").Value = 2529.6667
$ws.Range("M68").Value = -1173.5
$ws.Range("N68").Value = -4027.6667
$ws.Range("H71").Value = 2286.8
$ws.Range("I71").Value = 1922.5
$ws.Range("J71").Value = 2529.6667
$ws.Range("K71").Value = 9612.5
$ws.Range("L71").Value = 12648.3335
$ws.Range("M71").Value = -5868.5
$ws.Range("N71").Value = -20136.3335
$ws.Range("H113").Value = 1908.3334
$ws.Range("I113").Value = 1902.9412
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1902.9412
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 267.0588
$ws.Range("N113").Value = -6340
$ws.Range("H132").Value = 8029.081
$ws.Range("I132").Value = 2840.15
$ws.Range("J132").Value = 14133.706
$ws.Range("K132").Value = 8520.450000000001
$ws.Range("L132").Value = 42401.118
$ws.Range("M132").Value = -5990.450000000001
$ws.Range("N132").Value = -47461.118

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 539.6429000000001
$ws.Range("I107").Value = 478.85715
$ws.Range("J107").Value = 600.4286
$ws.Range("K107").Value = 1436.57145
$ws.Range("L107").Value = 1801.2858
$ws.Range("M107").Value = 483.4285500000001
